$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1031.6666
$ws.Range("I11").Value = 1031.6666
$ws.Range("K11").Value = 1031.6666
$ws.Range("M11").Value = -891.6666

$ws.Range("H17").Value = 1727
$ws.Range("J17").Value = 1727
$ws.Range("L17").Value = 5181
$ws.Range("N17").Value = -5517

$ws.Range("H32").Value = 400
$ws.Range("J32").Value = 400
$ws.Range("L32").Value = 400
$ws.Range("N32").Value = -1052

$ws.Range("H40").Value = 3297.9666
$ws.Range("I40").Value = 3790.6365
$ws.Range("K40").Value = 3790.6365
$ws.Range("M40").Value = -3615.6365

$ws.Range("H55").Value = 150
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H69").Value = 18319.875
$ws.Range("I69").Value = 11498.5
$ws.Range("J69").Value = 25141.25
$ws.Range("K69").Value = 34495.5
$ws.Range("L69").Value = 75423.75
$ws.Range("M69").Value = -33621.5
$ws.Range("N69").Value = -77171.75

$ws.Range("H72").Value = 18319.875
$ws.Range("I72").Value = 11498.5
$ws.Range("J72").Value = 25141.25
$ws.Range("K72").Value = 103486.5
$ws.Range("L72").Value = 226271.25
$ws.Range("M72").Value = -99118.5
$ws.Range("N72").Value = -235007.25

$ws.Range("H98").Value = 43479504
$ws.Range("I98").Value = 47619956
$ws.Range("K98").Value = 47619956
$ws.Range("M98").Value = -47618458

$ws.Range("H107").Value = 610.4761999999999
$ws.Range("I107").Value = 861.7857
$ws.Range("K107").Value = 861.7857
$ws.Range("M107").Value = 1058.2143

$ws.Range("H108").Value = 68000
$ws.Range("J108").Value = 68000
$ws.Range("L108").Value = 68000
$ws.Range("N108").Value = -75680

$ws.Range("H110").Value = 31500
$ws.Range("J110").Value = 31500
$ws.Range("L110").Value = 31500
$ws.Range("N110").Value = -39680

$ws.Range("H114").Value = 70000
$ws.Range("J114").Value = 70000
$ws.Range("L114").Value = 70000
$ws.Range("N114").Value = -78678

$ws.Range("H122").Value = 43479504
$ws.Range("I122").Value = 47619956
$ws.Range("K122").Value = 142859868
$ws.Range("M122").Value = -142857418

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10455393
$ws.Range("I32").Value = 17903080
$ws.Range("K32").Value = 17903080
$ws.Range("M32").Value = -17902793

$ws.Range("H132").Value = 6806412.5
$ws.Range("I132").Value = 9806152
$ws.Range("K132").Value = 29418456
$ws.Range("M132").Value = -29415926

$ws.Range("H137").Value = 15000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 99977
$ws.Range("J130").Value = 99977
$ws.Range("L130").Value = 99977
$ws.Range("N130").Value = -110017

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 45450
$ws.Range("I2").Value = 45450
$ws.Range("K2").Value = 45450
$ws.Range("M2").Value = -45337

$ws.Range("H20").Value = 133990
$ws.Range("J20").Value = 133990
$ws.Range("L20").Value = 133990
$ws.Range("N20").Value = -134462

$ws.Range("H30").Value = 133990
$ws.Range("J30").Value = 133990
$ws.Range("L30").Value = 133990
$ws.Range("N30").Value = -134172

$ws.Range("H31").Value = 1798333.8
$ws.Range("I31").Value = 2264
$ws.Range("K31").Value = 2264
$ws.Range("M31").Value = -1969

$ws.Range("H34").Value = 1798333.8
$ws.Range("I34").Value = 2264
$ws.Range("K34").Value = 2264
$ws.Range("M34").Value = -2062

$ws.Range("H58").Value = 7117.4287
$ws.Range("I58").Value = 5852.25
$ws.Range("K58").Value = 5852.25
$ws.Range("M58").Value = -5649.25

$ws.Range("H94").Value = 2892.6453
$ws.Range("I94").Value = 2309.6667
$ws.Range("J94").Value = 3439.1875
$ws.Range("K94").Value = 2309.6667
$ws.Range("L94").Value = 3439.1875
$ws.Range("M94").Value = -1858.6667
$ws.Range("N94").Value = -4341.1875

$ws.Range("H111").Value = 79930
$ws.Range("J111").Value = 79930
$ws.Range("L111").Value = 79930
$ws.Range("N111").Value = -88110

$ws.Range("H124").Value = 66220.336
$ws.Range("J124").Value = 66220.336
$ws.Range("L124").Value = 66220.336
$ws.Range("N124").Value = -71130.336

$ws.Range("H128").Value = 133990
$ws.Range("J128").Value = 133990
$ws.Range("L128").Value = 133990
$ws.Range("N128").Value = -143950

$ws.Range("H129").Value = 65995
$ws.Range("J129").Value = 65995
$ws.Range("L129").Value = 65995
$ws.Range("N129").Value = -75995

$ws.Range("H130").Value = 91999.75
$ws.Range("J130").Value = 91999.75
$ws.Range("L130").Value = 91999.75
$ws.Range("N130").Value = -102039.75

$ws.Range("H132").Value = 4505.077
$ws.Range("I132").Value = 1377.2858
$ws.Range("K132").Value = 4131.857400000001
$ws.Range("M132").Value = -1601.857400000001

$ws.Range("H134").Value = 3756.8215
$ws.Range("I134").Value = 1533.8
$ws.Range("K134").Value = 4601.4
$ws.Range("M134").Value = -2066.4

$ws.Range("H136").Value = 7117.4287
$ws.Range("I136").Value = 5852.25
$ws.Range("K136").Value = 17556.75
$ws.Range("M136").Value = -15006.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9515707
$ws.Range("J4").Value = 14344405
$ws.Range("L4").Value = 43033215
$ws.Range("N4").Value = -43033439

$ws.Range("H9").Value = 5895
$ws.Range("I9").Value = 10393.333
$ws.Range("J9").Value = 1396.6666
$ws.Range("K9").Value = 31179.999
$ws.Range("L9").Value = 4189.9998
$ws.Range("M9").Value = -30955.999
$ws.Range("N9").Value = -4637.9998

$ws.Range("H141").Value = 342775.44
$ws.Range("J141").Value = 15000
$ws.Range("L141").Value = 45000
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6251.4287
$ws.Range("I70").Value = 5715.25
$ws.Range("J70").Value = 6966.3335
$ws.Range("K70").Value = 5715.25
$ws.Range("L70").Value = 6966.3335
$ws.Range("M70").Value = -5445.25
$ws.Range("N70").Value = -7506.3335

$ws.Range("H73").Value = 6251.4287
$ws.Range("I73").Value = 5715.25
$ws.Range("J73").Value = 6966.3335
$ws.Range("K73").Value = 5715.25
$ws.Range("L73").Value = 6966.3335
$ws.Range("M73").Value = -4779.25
$ws.Range("N73").Value = -8838.333500000001

$ws.Range("H108").Value = 119979.336
$ws.Range("J108").Value = 119979.336
$ws.Range("L108").Value = 119979.336
$ws.Range("N108").Value = -127659.336

$ws.Range("H113").Value = 4456.5264
$ws.Range("I113").Value = 3770.889
$ws.Range("K113").Value = 3770.889
$ws.Range("M113").Value = -1600.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1202.9445
$ws.Range("I22").Value = 1092.1818
$ws.Range("K22").Value = 1092.1818
$ws.Range("M22").Value = -797.1818000000001

$ws.Range("H27").Value = 1202.9445
$ws.Range("I27").Value = 1092.1818
$ws.Range("K27").Value = 1092.1818
$ws.Range("M27").Value = -985.1818000000001

$ws.Range("H46").Value = 3713.7778
$ws.Range("I46").Value = 2841.9
$ws.Range("J46").Value = 4803.625
$ws.Range("K46").Value = 2841.9
$ws.Range("L46").Value = 4803.625
$ws.Range("M46").Value = -2653.9
$ws.Range("N46").Value = -5179.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1525.8334
$ws.Range("I113").Value = 1320.2858
$ws.Range("K113").Value = 3960.8574
$ws.Range("M113").Value = -1790.8574

$ws.Range("H132").Value = 3092616
$ws.Range("I132").Value = 5829.9
$ws.Range("K132").Value = 17489.7
$ws.Range("M132").Value = -14959.7
